# Reset the "Weight" (column E) values that are not already 1 back to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(3, 6, 7, 10, 11, 13, 14, 16, 17, 18, 23, 24, 25, 26, 27, 28, 29, 33, 34, 35, 36, 37, 38, 40, 41, 44)

foreach ($r in $rows) {
    $ws.Range("E$r").Value = 1
}
